$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-7: numeric data in columns A:C
$data = @(
    @(88, 46, 57),
    @(89, 38, 12),
    @(23, 59, 78),
    @(56, 21, 98),
    @(24, 18, 43),
    @(34, 15, 67)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Row 8: header-like text row across A:D
$ws.Cells.Item(8, 1).Value = "type1"
$ws.Cells.Item(8, 2).Value = "type2"
$ws.Cells.Item(8, 3).Value = "type3"
$ws.Cells.Item(8, 4).Value = "type4"

# Row 9: another text row across A:D
$ws.Cells.Item(9, 1).Value = "tup1"
$ws.Cells.Item(9, 2).Value = "tup2"
$ws.Cells.Item(9, 3).Value = "tup3"
$ws.Cells.Item(9, 4).Value = "tup4"
